$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Count (Total: 223)"
$ws.Range("B2").Value = 135
$ws.Range("B3").Value = 73
$ws.Range("B4").Value = 15
